$d = $word.ActiveDocument

$replacements = @(
    @("561÷6=93, 3", "319÷9=35, 4"),
    @("174÷7=24, 6", "487÷9=54, 1"),
    @("544÷9=60, 4", "516÷5=103, 1"),
    @("493÷3=164, 1", "698÷2=349, 0"),
    @("562÷6=93, 4", "200÷7=28, 4"),
    @("195÷4=48, 3", "836÷9=92, 8"),
    @("729÷8=91, 1", "360÷9=40, 0"),
    @("783÷4=195, 3", "694÷8=86, 6"),
    @("777÷6=129, 3", "711÷8=88, 7"),
    @("101÷2=50, 1", "417÷5=83, 2"),
    @("942÷9=104, 6", "390÷9=43, 3"),
    @("426÷5=85, 1", "382÷2=191, 0"),
    @("349÷5=69, 4", "331÷2=165, 1"),
    @("861÷4=215, 1", "571÷5=114, 1"),
    @("584÷4=146, 0", "595÷9=66, 1"),
    @("309÷5=61, 4", "131÷2=65, 1"),
    @("621÷2=310, 1", "433÷7=61, 6"),
    @("237÷8=29, 5", "676÷8=84, 4"),
    @("469÷4=117, 1", "566÷7=80, 6"),
    @("465÷2=232, 1", "725÷9=80, 5"),
    @("379÷3=126, 1", "203÷4=50, 3"),
    @("598÷5=119, 3", "287÷6=47, 5"),
    @("285÷7=40, 5", "616÷8=77, 0"),
    @("116÷2=58, 0", "482÷2=241, 0"),
    @("305÷2=152, 1", "657÷7=93, 6")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
